$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing requirement texts (renumbering SSS-000N -> SSS-00N, etc.) ---
$ws.Cells.Item(3,1).Value  = "SSS-001  O sistema DEVE receber a solicitação de serviço"
$ws.Cells.Item(4,1).Value  = "SSS-002 O sistema DEVE armazenar os dados do cliente e veículo"
$ws.Cells.Item(5,1).Value  = "SSS-003 O sistema DEVE fornecer orçamento ao cliente "
$ws.Cells.Item(6,1).Value  = "SSS-002 O sistema DEVE permitir que o Cliente acompanhe as situações de status de seu veículo"
$ws.Cells.Item(7,1).Value  = "SSS-004 O sistema DEVE permitir que o estoquista registre o material"
$ws.Cells.Item(8,1).Value  = "SSS-005 O sistema DEVE informar ao estoquista o material  existente "
$ws.Cells.Item(9,1).Value  = "SSS-006 O sistema DEVE informar ao estoquista os materiais inexistente"
$ws.Cells.Item(10,1).Value = "SSS-007 O sistema DEVE informar ao avaliador o cálculo de material para o serviço"
$ws.Cells.Item(11,1).Value = "SSS-008 O sistema DEVE registrar entrada de lucros e despesas da oficina"
$ws.Cells.Item(12,1).Value = "SSS-009 O sistema DEVE permitir que o Setor Financeiro consulte os lucros e despesas"

# Rows 13-14 keep their numbering (SSS-0010 / SSS-0011), text unchanged
$ws.Cells.Item(13,1).Value = "SSS-0010 O sistema DEVE registrar todos os funcionários, cada um com a sua digital para controle de ponto  "
$ws.Cells.Item(14,1).Value = "SSS-0011 O sistema DEVE oferecer ao supervisor a opção de consulta do controle de ponto dos funcionários"

# Row 15 (SSS-0012) gains an extra clause
$ws.Cells.Item(15,1).Value = "SSS-0012 O sistema DEVE consultar os dados do veículo para confirmar se não consta como roubado/furtado antes de executar qualquer ação, migrando para o site do Detran"

# --- Append brand-new requirements SSS-0013 .. SSS-0025 in rows 17-29 ---
$ws.Cells.Item(17,1).Value = "SSS-0013 O sistema DEVE gerar relatórios dinâmicos com gráficos baseado na quantidade de consertos finalizados"
$ws.Cells.Item(18,1).Value = "SSS-0014 O sistema DEVE gerar relatórios dinâmicos com gráficos baseado na quantidade de consertos pendentes"
$ws.Cells.Item(19,1).Value = "SSS-0015 O sistema DEVE gerar  relatórios dinâmicos com gráficos baseado na quantidade de clientes que pedem orçamento"
$ws.Cells.Item(20,1).Value = "SSS-0016 O sistema DEVE gerar  relatórios dinâmicos com gráficos baseado em qual serviço é executado com mais propriedade"
$ws.Cells.Item(21,1).Value = "SSS-0017 O sistema DEVE gerar  relatórios dinâmicos com gráficos baseado em quais peças possuem maior saída "
$ws.Cells.Item(22,1).Value = "SSS-0018 O sistema DEVE gerar  relatórios dinâmicos com gráficos baseado em quais peças possuem menos saída "
$ws.Cells.Item(23,1).Value = "SSS-0019 O sistema DEVE gerar  relatórios dinâmicos com gráficos baseado na produtividade de cada funcionário"
$ws.Cells.Item(24,1).Value = "SSS-0020 O sistema DEVE informar se o cliente já consta no banco de dados, se solicitou algum  serviço "
$ws.Cells.Item(25,1).Value = "SSS-0021 O sistema DEVE registrar nos dias de trabalho os horários de entrada, pausa e saída de seus funcionário"
$ws.Cells.Item(26,1).Value = "SSS-0022 O sistema DEVE gerar os lucros e depesas de fechamento de cada mês "
$ws.Cells.Item(27,1).Value = "SSS-0023 O sistema DEVE calcular o prazo de conserto "
$ws.Cells.Item(28,1).Value = "SSS-0024 O sistema DEVE gerar um acesso exclusivo para cada cliente"
$ws.Cells.Item(29,1).Value = "SSS-0025 O sistema DEVE registar agendamentos com datas e horários para não ocorrer enganos "

# Copy the existing "left align" formatting (used by rows 3-12) onto the new rows 17-29
$src = $ws.Range("A12")
$dest = $ws.Range("A17:A29")
$src.Copy()
$dest.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Widen column A to fit the new, longer text (closest achievable value to 124.42578125)
$ws.Columns.Item(1).ColumnWidth = 123.6

# Move the active selection like the author left it
$ws.Range("A9").Select()
